# CoronaVirus US Lockdown Forecast - update the forecast model.
#
# The forecast column (I16:I28) previously projected I16 itself from I15
# via the shared growth formula. The actual reported value for 4/6 (row 16)
# is now known, so I16 becomes a hard-coded actual value (matching the
# "actuals" look/format of I15), and the projection formula now starts
# fresh from I17, growing off of the new I16 actual.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give I16 the same "actual data" formatting as I15 (fill/number format),
# replacing its previous "forecast" formatting, then overwrite it with the
# real reported value.
$ws.Range("I15").Copy()
$ws.Range("I16").PasteSpecial(-4122)
$ws.Range("I16").Value = 33546

# Re-anchor the forecast formula so it now starts at I17, growing off of
# I16 (the new actual), propagating down through I28. Writing the same
# formula text across the whole range lets relative references (I and M)
# adjust per row, same as Excel's own fill-down.
$ws.Range("I17:I28").Formula = "=I16*(1+AVERAGE(M12:M16))"

# Reflect where the user's cursor ended up after the edit.
$ws.Range("I17").Select()
